$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the server data row (row 2) - order matches shared-string insertion order
$ws.Range("F2").Value = "127.0.0.1"

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "ProxyServer_1"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "ProxyServer_1"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "000105001"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 5001

# Update the selection to match the author's final cursor position
$ws.Range("G4").Select()
